$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Merge the two adjacent runs "(from the database at earlier stages)" and
# " " (a trailing space, both red/FF0000) into a single run with the
# combined text "(from the database at earlier stages) ". A find/replace
# whose search and replacement text are identical (but span both runs)
# makes Word normalise the run boundaries, collapsing the two runs into
# one that keeps the first run's formatting.
$rng1 = $d.Content
$rng1.Find.Execute( `
    "(from the database at earlier stages) ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "(from the database at earlier stages) ", 2)

# --- Change 2 ---------------------------------------------------------
# Append "b." right after "The accepted image size is 50 M" so the
# sentence reads "...50 Mb.", matching the run's existing character
# formatting (10pt / complex-script 10pt / white highlight).
$rng2 = $d.Content
$rng2.Find.Execute( `
    "The accepted image size is 50 M", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.Select()
$sel = $word.Selection
$sel.InsertAfter("b.")
$sel.HighlightColorIndex = 8
$sel.Font.Size = 10
$sel.Font.SizeBi = 10
